$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the obsolete "Added instantiation of ... BC emissions" / Review /
#    In Progress entry that used to live in row 99. This shifts the old rows
#    100-102 up to 99-101.
# ---------------------------------------------------------------------------
$ws.Rows.Item(99).Delete()

# ---------------------------------------------------------------------------
# 2) Fix up the (now shifted) rows 99-101: update change numbers, statuses
#    and the proposed/review/committed dates.
# ---------------------------------------------------------------------------

# Row 99: Created C1.2.add_NC_emissions_EDGAR.R ...
$ws.Range("B99").Value = 95
$ws.Range("D99").Value = "Committed"
$ws.Range("E99").Value = 42377
$ws.Range("F99").Value = 42377
$ws.Range("G99").Value = 42377

# Row 100: Created NC_EDGAR_sector_mapping.csv ...
$ws.Range("B100").Value = 96
$ws.Range("D100").Value = "Committed"
$ws.Range("E100").Value = 42380
$ws.Range("E100").Copy()
$ws.Range("F100:G100").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F100").Value = 42016
$ws.Range("G100").Value = 42017
$ws.Range("H100").Value = "5edb7be"

# Row 101: Added EDGARcheck function to analysis_functions.R
$ws.Range("B101").Value = 97
$ws.Range("D101").Value = "Committed"
$ws.Range("E101").Value = 42380
$ws.Range("E101").Copy()
$ws.Range("F101:G101").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F101").Value = 42016
$ws.Range("G101").Value = 42017
$ws.Range("H101").Value = "5edb7be"

# ---------------------------------------------------------------------------
# 3) Append the two new rows describing the IO_functions.R / UNFCCC work.
# ---------------------------------------------------------------------------

# Row 102: Upgraded IO_functions.R readData function ...
$ws.Range("A102").Value = "Upgraded IO_functions.R readData function to include ability to read one, all, or a select list of .csv files from within a .zip file. Added listZippedFiles function."
$ws.Range("B102").Value = 98
$ws.Range("C102").Value = "Jon Seibert"
$ws.Range("D102").Value = "Committed"
$ws.Range("E101").Copy()
$ws.Range("E102:G102").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E102").Value = 42384
$ws.Range("F102").Value = 42387
$ws.Range("G102").Value = 42387
$ws.Range("H102").Value = "ce6f6a3"

$ws.Range("A102:H102").RowHeight = 54.75

# Row 103: Renamed E.UNFCCC_SO2_emissions.R to E.UNFCCC_emissions.R ...
$ws.Range("A103").Value = "Renamed E.UNFCCC_SO2_emissions.R to E.UNFCCC_emissions.R, updated to use new readData .zip features to read all data from within large .zip files, added dummy output for species without present input data."
$ws.Range("B103").Value = 99
$ws.Range("C103").Value = "Jon Seibert"
$ws.Range("D103").Value = "Committed"
$ws.Range("E102").Copy()
$ws.Range("E103:G103").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E103").Value = 42384
$ws.Range("F103").Value = 42387
$ws.Range("G103").Value = 42387
$ws.Range("H103").Value = "ce6f6a3"

$ws.Range("A103:H103").RowHeight = 73.5

# ---------------------------------------------------------------------------
# 4) Set the row heights that resulted from the edit (some already match
#    after the delete-shift above, but set them explicitly to be safe).
# ---------------------------------------------------------------------------
$ws.Rows.Item(99).RowHeight = 64.5
$ws.Rows.Item(100).RowHeight = 62.25
$ws.Rows.Item(101).RowHeight = 32.25
$ws.Rows.Item(102).RowHeight = 54.75
$ws.Rows.Item(103).RowHeight = 73.5

# ---------------------------------------------------------------------------
# 5) Trailing empty row left behind by the author.
# ---------------------------------------------------------------------------
$ws.Rows.Item(104).RowHeight = 36.75

# ---------------------------------------------------------------------------
# 6) Selection, as left by the author after editing.
# ---------------------------------------------------------------------------
$ws.Range("C105").Select()
